$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.535914897918701
$ws.Range("B1").Value = 2.637543916702271
$ws.Range("C1").Value = 1.811471343040466
$ws.Range("D1").Value = 1.61442756652832
$ws.Range("E1").Value = 1.54661762714386
